$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 9492
$ws.Range("I34").Value = 9492
$ws.Range("K34").Value = 9492
$ws.Range("M34").Value = -9289
$ws.Range("H36").Value = 9492
$ws.Range("I36").Value = 9492
$ws.Range("K36").Value = 9492
$ws.Range("M36").Value = -8777
$ws.Range("H47").Value = 17533
$ws.Range("J47").Value = 17799.5
$ws.Range("L47").Value = 17799.5
$ws.Range("N47").Value = -19743.5
$ws.Range("H111").Value = 1543.4
$ws.Range("I111").Value = 1479.25
$ws.Range("K111").Value = 4437.75
$ws.Range("M111").Value = -1370.75
$ws.Range("H132").Value = 3009
$ws.Range("I132").Value = 679
$ws.Range("K132").Value = 2037
$ws.Range("M132").Value = 493
$ws.Range("H137").Value = 1765.8667
$ws.Range("I137").Value = 1612.125
$ws.Range("K137").Value = 4836.375
$ws.Range("M137").Value = -2286.375
$ws.Range("H138").Value = 2100.16
$ws.Range("J138").Value = 2081.5715
$ws.Range("L138").Value = 6244.7145
$ws.Range("N138").Value = -16524.7145
$ws.Range("H141").Value = 1191.1111
$ws.Range("I141").Value = 1191.1111
$ws.Range("K141").Value = 3573.3333
$ws.Range("M141").Value = 1606.6667

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1597.6923
$ws.Range("I74").Value = 1622.5
$ws.Range("J74").Value = 1300
$ws.Range("K74").Value = 1622.5
$ws.Range("L74").Value = 1300
$ws.Range("M74").Value = -748.5
$ws.Range("N74").Value = -3048
$ws.Range("H77").Value = 1597.6923
$ws.Range("I77").Value = 1622.5
$ws.Range("J77").Value = 1300
$ws.Range("K77").Value = 8112.5
$ws.Range("L77").Value = 6500
$ws.Range("M77").Value = -3744.5
$ws.Range("N77").Value = -15236

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5752.3076
$ws.Range("I107").Value = 4397.875
$ws.Range("K107").Value = 4397.875
$ws.Range("M107").Value = -2477.875
$ws.Range("H134").Value = 1480.7646
$ws.Range("I134").Value = 1260.875
$ws.Range("K134").Value = 3782.625
$ws.Range("M134").Value = -1247.625

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 86.63158
$ws.Range("I7").Value = 47.545456
$ws.Range("J7").Value = 140.375
$ws.Range("K7").Value = 47.545456
$ws.Range("L7").Value = 140.375
$ws.Range("M7").Value = 65.454544
$ws.Range("N7").Value = -366.375
$ws.Range("H22").Value = 992.06665
$ws.Range("I22").Value = 658.1111
$ws.Range("K22").Value = 658.1111
$ws.Range("M22").Value = -308.1111
$ws.Range("H31").Value = 2432.7
$ws.Range("I31").Value = 1920.2858
$ws.Range("J31").Value = 3628.3333
$ws.Range("K31").Value = 1920.2858
$ws.Range("L31").Value = 3628.3333
$ws.Range("M31").Value = -1625.2858
$ws.Range("N31").Value = -4218.3333
$ws.Range("H34").Value = 2432.7
$ws.Range("I34").Value = 1920.2858
$ws.Range("J34").Value = 3628.3333
$ws.Range("K34").Value = 1920.2858
$ws.Range("L34").Value = 3628.3333
$ws.Range("M34").Value = -1718.2858
$ws.Range("N34").Value = -4032.3333
$ws.Range("H94").Value = 170147.14
$ws.Range("I94").Value = 196002.67
$ws.Range("K94").Value = 196002.67
$ws.Range("M94").Value = -195551.67
$ws.Range("H122").Value = 1914
$ws.Range("I122").Value = 2083.1667
$ws.Range("K122").Value = 6249.500100000001
$ws.Range("M122").Value = -3799.500100000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 130107.82
$ws.Range("I2").Value = 91706.414
$ws.Range("J2").Value = 222271.2
$ws.Range("K2").Value = 550238.4840000001
$ws.Range("L2").Value = 1333627.2
$ws.Range("M2").Value = -550125.4840000001
$ws.Range("N2").Value = -1333853.2
$ws.Range("H37").Value = 59950
$ws.Range("J37").Value = 59950
$ws.Range("L37").Value = 179850
$ws.Range("N37").Value = -180074
$ws.Range("H55").Value = 1024.5
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354
$ws.Range("H64").Value = 3366
$ws.Range("I64").Value = 549.5
$ws.Range("K64").Value = 1648.5
$ws.Range("M64").Value = -1378.5
$ws.Range("H67").Value = 3366
$ws.Range("I67").Value = 549.5
$ws.Range("K67").Value = 1648.5
$ws.Range("M67").Value = -712.5
$ws.Range("H116").Value = 1500
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1058
$ws.Range("N116").ClearContents()
$ws.Range("H140").Value = 7957.731
$ws.Range("I140").Value = 1244.125
$ws.Range("J140").Value = 18699.5
$ws.Range("K140").Value = 3732.375
$ws.Range("L140").Value = 56098.5
$ws.Range("M140").Value = 1447.625
$ws.Range("N140").Value = -66458.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1204.1428
$ws.Range("I22").Value = 1087.1818
$ws.Range("J22").Value = 1633
$ws.Range("K22").Value = 1087.1818
$ws.Range("L22").Value = 1633
$ws.Range("M22").Value = -792.1818000000001
$ws.Range("N22").Value = -2223
$ws.Range("H27").Value = 1204.1428
$ws.Range("I27").Value = 1087.1818
$ws.Range("J27").Value = 1633
$ws.Range("K27").Value = 1087.1818
$ws.Range("L27").Value = 1633
$ws.Range("M27").Value = -980.1818000000001
$ws.Range("N27").Value = -1847
$ws.Range("H82").Value = 1439.3158
$ws.Range("I82").Value = 1436.6
$ws.Range("J82").Value = 1449.5
$ws.Range("K82").Value = 1436.6
$ws.Range("L82").Value = 1449.5
$ws.Range("M82").Value = -1075.6
$ws.Range("N82").Value = -2171.5
$ws.Range("H85").Value = 1439.3158
$ws.Range("I85").Value = 1436.6
$ws.Range("J85").Value = 1449.5
$ws.Range("K85").Value = 1436.6
$ws.Range("L85").Value = 1449.5
$ws.Range("M85").Value = -188.5999999999999
$ws.Range("N85").Value = -3945.5
$ws.Range("H132").Value = 1936.9474
$ws.Range("I132").Value = 1766.7778
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5300.3334
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2770.3334
$ws.Range("N132").Value = -20060

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 44999
$ws.Range("J49").Value = 44999
$ws.Range("L49").Value = 44999
$ws.Range("N49").Value = -45459
$ws.Range("H126").Value = 2608.5833
$ws.Range("I126").Value = 2630.75
$ws.Range("K126").Value = 7892.25
$ws.Range("M126").Value = -5422.25
$ws.Range("H132").Value = 2125.5186
$ws.Range("I132").Value = 2273.56
$ws.Range("J132").Value = 275
$ws.Range("K132").Value = 6820.68
$ws.Range("L132").Value = 825
$ws.Range("M132").Value = -4290.68
$ws.Range("N132").Value = -5885
$ws.Range("H136").Value = 2840.8518
$ws.Range("I136").Value = 2843.739
$ws.Range("J136").Value = 2824.25
$ws.Range("K136").Value = 8531.217000000001
$ws.Range("L136").Value = 8472.75
$ws.Range("M136").Value = -5981.217000000001
$ws.Range("N136").Value = -13572.75
